$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TranslationRow {
    param($Row, $Key, $Translation)
    $addr = "A" + $Row + ":C" + $Row
    $target = $ws.Range($addr)
    $target.WrapText = $true
    $target.Font.Size = 10
    $ws.Range("A" + $Row).Value = "cs"
    $ws.Range("B" + $Row).Value = $Key
    $ws.Range("C" + $Row).Value = $Translation
}

# New "public 404 / sign-up / registration" rows, written in the same
# order they were authored so the shared-string table lines up with the
# source workbook.
Set-TranslationRow 18 "public.sign-up.menu" "Registrace"
Set-TranslationRow 19 "public.404.title" "Tady nic není!"
Set-TranslationRow 20 "public.404.title" "I těm nejlepším se to stane…"
Set-TranslationRow 21 "public.404.subtitle" "Bohužel jste narazili na stránku, která neexistuje. Je to divné, ale je to tak."
Set-TranslationRow 22 "public.404.back" "Zpět"
Set-TranslationRow 23 "public.404.home" "Domů"
Set-TranslationRow 24 "public.development-notice.alert" "Aplikace je stále ve vývoji a mnoho vychytávek chybí, nicméně je aktivně vyvíjena i používána."

# Existing row 6 key gets renamed as part of the same change
$ws.Range("B6").Value = "public.sign-in.title"

Set-TranslationRow 25 "public.sign-up.title" "Registrace"
Set-TranslationRow 26 "public.sign-up.subtitle" "Po registraci získáte okamžitý přístup do aplikace a můžete začít zkoumat zajímavé možnosti, kterými disponuje."
Set-TranslationRow 27 "user.name.label" "Vaše jméno"
Set-TranslationRow 28 "user.name.label.tooltip" "Jedná se víceméně o jakékoli jméno, kterým se chcete prezentovat. Pro přihlášení bude použit Váš email."
Set-TranslationRow 29 "user.password2.label" "Kontrola hesla"
Set-TranslationRow 30 "public.sign-up.form.submit.label" "Registrovat se"
Set-TranslationRow 31 "user.email.label" "E-mail"
Set-TranslationRow 32 "user.email.label.tooltip" "Emailová adresa musí být v systému unikátní a slouží dále pro přihlášení do systému."
Set-TranslationRow 33 "user.password.label.required" "Bez hesla to bohužel nejde, vyplňte jej prosím."
Set-TranslationRow 34 "user.password2.label.required" "Kontrolní heslo je vyžadováno."
Set-TranslationRow 35 "user.name.label.required" "Vaše jméno je vyžadováno, vyplňte jej prosím."
Set-TranslationRow 36 "user.email.label.required" "E-mailová adresa slouží k přihlášení do systému, je povinná."
Set-TranslationRow 37 "user.password2.label.mismatch" "Hesla se neshoduji!"

# A few rows wrap onto a second line in the original workbook.
$ws.Rows.Item(24).RowHeight = 26.25
$ws.Rows.Item(26).RowHeight = 26.25
$ws.Rows.Item(28).RowHeight = 26.25

$ws.Range("B32").Select()
